$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 89

# Date-like and pure-integer-like text must be forced to Text format first,
# otherwise Excel auto-coerces "2023-06-29" into a date serial and "26" into
# a number. ClearFormats() afterwards drops the residual Text numberformat
# style so the cell keeps the default (no explicit style), matching the
# other text cells in the sheet, while the stored value stays a string.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2023-06-29"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "08:37:45"

$ws.Cells.Item($row, 3).Value = "Thursday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "26"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 123079
$ws.Cells.Item($row, 6).Value = 134475
$ws.Cells.Item($row, 7).Value = 163581
$ws.Cells.Item($row, 8).Value = 134098
$ws.Cells.Item($row, 9).Value = 177212
$ws.Cells.Item($row, 10).Value = 115019
$ws.Cells.Item($row, 11).Value = 204255
$ws.Cells.Item($row, 12).Value = 226326
$ws.Cells.Item($row, 13).Value = 176371
$ws.Cells.Item($row, 14).Value = 104497
$ws.Cells.Item($row, 15).Value = 39740
$ws.Cells.Item($row, 16).Value = 33717
$ws.Cells.Item($row, 17).Value = 52446
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 35620
$ws.Cells.Item($row, 20).Value = -1
